$d = $word.ActiveDocument

# Locate the paragraph that follows the "Requisitos" entry we keep
# ("LOQ4064: ...") and remove the trailing boilerplate block that was
# appended by the site generator: a blank paragraph, the "Ver no
# Jupiter ..." line, and the "(c) 2020 ..." footer line.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "LOQ4064") {
        $target = $i
        break
    }
}

if ($target -ne $null) {
    $startPara = $d.Paragraphs.Item($target + 1)
    $endPara = $d.Paragraphs.Item($target + 3)
    $r = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $r.Delete()
}
